$wb = $excel.ActiveWorkbook

# --- Update the selection on the original sheet (听说教程3) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D2").Select()

# --- Add the new worksheet "视听说教程3" after the existing sheet ---
# The engine assigns sheetId based on (current sheet count + 1) at creation time,
# and reuses ids after a delete. The target file expects sheetId="3" for the new
# sheet even though the workbook only ends up with 2 sheets, so we create a
# throwaway sheet first to "use up" id 2, then create the real sheet (id 3),
# then remove the throwaway again.
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$tempSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$tempSheet.Name = "TempPlaceholder"

$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last2)
$newSheet.Name = "视听说教程3"

$toDelete = $wb.Worksheets.Item("TempPlaceholder")
$toDelete.Delete()

$target = $wb.Worksheets.Item("视听说教程3")
$target.Activate()

# --- Fill in the statistics table ---
# Values are written in the same order the original author entered them so
# that new shared-string entries land on the same indices as the source file.
$target.Range("A3").Value = "Unit"
$target.Range("B3").Value = "Ouside view"
$target.Range("B7").Value = "Listening in"
$target.Range("C3").Value = "Activity 1"
$target.Range("C4").Value = "Activity 2"
$target.Range("C5").Value = "Activity 3"
$target.Range("C6").Value = "Activity 4"
$target.Range("C7").Value = "Passage 1: Activity 1"
$target.Range("C8").Value = "Passage 1: Activity 2"
$target.Range("C9").Value = "Passage 1: Activity 3"
$target.Range("C10").Value = "Passage 2: Activity 1"
$target.Range("C11").Value = "Passage 2: Activity 2"
$target.Range("C12").Value = "Passage 2: Activity 3"
$target.Range("D2").Value = "Target: 10.13.54.81/book/book41/"
$target.Range("D3").Value = "dj31drag.php"
$target.Range("D4").Value = "dj34mc.php"
$target.Range("D5").Value = "dj31checkboxTable.php"
$target.Range("D6").Value = "dj34drag.php"
$target.Range("D7").Value = "dj44mc.php"
$target.Range("D8").Value = "dj41blank.php"
$target.Range("D12").Value = "dj42drag.php"
$target.Range("D10").Value = "dj42drag.php/dj44dragOne.php"
$target.Range("D11").Value = "dj45checkboxTable.php"

# --- Merge the Unit / section label cells ---
$target.Range("A3:A12").Merge()
$target.Range("B3:B6").Merge()
$target.Range("B7:B12").Merge()

# --- Center alignment for the labeled header / merged cells ---
$target.Range("D2").HorizontalAlignment = -4108
$target.Range("D2").VerticalAlignment = -4108
$target.Range("A3:A12").HorizontalAlignment = -4108
$target.Range("A3:A12").VerticalAlignment = -4108
$target.Range("B3:B12").HorizontalAlignment = -4108
$target.Range("B3:B12").VerticalAlignment = -4108

# --- Column widths approximating the source's best-fit widths ---
$target.Columns.Item(1).ColumnWidth = 4.75
$target.Columns.Item(2).ColumnWidth = 11.375
$target.Columns.Item(3).ColumnWidth = 18.75
$target.Columns.Item(4).ColumnWidth = 31.75

# --- Final selection on the new sheet ---
$target.Range("F12").Select()
